$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")

# Insert a new row at row 2, shifting existing rows 2..135 down to 3..136
$ws.Rows("2:2").Insert()

# Fill in the new row 2 with the new transaction data
$ws.Range("E2").Value = "Withdrawal"
$ws.Range("N2").Value = "Credit Card"
$ws.Range("P2").Value = "Tradeprof"
$ws.Range("T2").Value = 269.29750000000001

# Update the view: top-left cell and selection
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("F15").Select()
